$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row -> new B value (text). $null means B is unchanged (row 38).
$updates = [ordered]@{
    2  = "6.69"
    8  = "8.50"
    14 = "2.99"
    20 = "12.85"
    26 = "11.32"
    32 = "27.62"
    38 = $null
    44 = "11.15"
    50 = "11.49"
    56 = "31.82"
    62 = "11.71"
    68 = "13.03"
    74 = "16.66"
}

foreach ($row in $updates.Keys) {
    Set-TextValue "A$row" "2025/12/30"
    $newB = $updates[$row]
    if ($newB -ne $null) {
        Set-TextValue "B$row" $newB
    }
}
